# "Magic sheet Names moved to ExcelProperties class"
# The sheet literal "REFERENCES" is replaced by the new constant "STOR_LOCKS".
# Renaming the sheet also keeps any defined names / autofilter references
# that point at it (e.g. the hidden _FilterDatabase name) in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REFERENCES")
$ws.Name = "STOR_LOCKS"
